$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.353.66"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.881.39"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'0.7113"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'242.81"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.08028"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.07%  "

$ws.Range("D9").Value = "'0.3139"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("D11").Value = "'0.08328"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").Value = "1.934.02"
$ws.Range("E12").Value = "  +3.43%  "

$ws.Range("D13").Value = "'5.265"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").Value = "'94.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.00%  "

$ws.Range("D15").Value = "'0.7189"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "'6.366"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.37%  "

$ws.Range("D17").Value = "'0.000008708"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.74%  "

$ws.Range("D18").Value = "29.401.71"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").Value = "'243.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").Value = "'13.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "2.121.21"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D23").Value = "'7.840"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  -1.87%  "

$ws.Range("D26").Value = "'163.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'9.094"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").Value = "'18.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = "'4.360"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.29%  "

$ws.Range("E32").Value = "  -6.39%  "

$ws.Range("D33").Value = "'0.05395"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("D34").Value = "'1.943"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'0.7810"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.93%  "

$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("D37").Value = "'2.688"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "'0.01887"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").Value = "1.268.95"
$ws.Range("E39").Value = "  +4.45%  "

$ws.Range("D40").Value = "'2.747"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("D41").Value = "'6.543"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").Value = "'0.9202"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.64%  "

$ws.Range("D43").Value = "'113.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.89%  "

$ws.Range("D44").Value = "'74.77"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.68%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("E46").Value = "  +3.97%  "

$ws.Range("D47").Value = "2.037.36"
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").Value = "'0.5223"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").Value = "'9.563"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").Value = "'0.4388"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.52%  "
